$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> [Price(D), Volume(E)] updated values. $null means D is unchanged.
$updates = @{
    2  = @("27.561.79", "  -0.13%  ")
    3  = @("1.648.68",  "  -0.53%  ")
    4  = @($null,       "  +0.04%  ")
    5  = @("212.76",    "  -0.99%  ")
    6  = @($null,       "  +5.33%  ")
    7  = @($null,       "  +0.03%  ")
    8  = @("23.57",     "  -2.47%  ")
    9  = @($null,       "  -1.47%  ")
    10 = @($null,       "  -1.20%  ")
    11 = @("0.0893",    "  +1.53%  ")
    12 = @("1.882.35",  "  -0.54%  ")
    13 = @("1.642.84",  "  -0.99%  ")
    14 = @("0.587",     "  +3.82%  ")
    15 = @($null,       "  -2.12%  ")
    16 = @("64.56",     "  -2.03%  ")
    17 = @("27.542.46", "  -0.07%  ")
    18 = @("231.48",    "  -3.91%  ")
    19 = @($null,       "  -0.36%  ")
    20 = @("7.54",      "  -0.90%  ")
    21 = @($null,       "  +0.08%  ")
    22 = @($null,       "  -3.28%  ")
    23 = @("9.78",      "  +4.39%  ")
    24 = @($null,       "  -2.25%  ")
    25 = @("148.82",    "  +1.96%  ")
    26 = @($null,       "  -2.47%  ")
    27 = @($null,       "  +1.51%  ")
    28 = @($null,       "  +0.19%  ")
    29 = @("15.62",     "  -4.24%  ")
    30 = @($null,       "  -2.29%  ")
    31 = @("0.0487",    "  -3.02%  ")
    32 = @("3.31",      "  -0.71%  ")
    33 = @("3.21",      "  +2.87%  ")
    34 = @("1.433.16",  "  -1.83%  ")
    35 = @("1.60",      "  +1.71%  ")
    36 = @($null,       "  -0.29%  ")
    37 = @($null,       "  -0.70%  ")
    38 = @("0.886",     "  -4.37%  ")
    39 = @("0.0167",    "  -2.99%  ")
    41 = @($null,       "  +0.07%  ")
    42 = @("0.820",     "  +3.90%  ")
    43 = @($null,       "  +2.24%  ")
    44 = @($null,       "  +1.55%  ")
    45 = @("65.40",     "  -6.11%  ")
    46 = @("1.791.94",  "  -0.48%  ")
    47 = @($null,       "  -1.25%  ")
    48 = @("88.09",     "  -0.43%  ")
    49 = @($null,       "  +0.26%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $price = $vals[0]
    $volume = $vals[1]
    if ($null -ne $price) {
        $ws.Range("D$row").Value = $price
    }
    $ws.Range("E$row").Value = $volume
}

# Rows 50 and 51 swap coin identity (Algorand <-> EnergySwap) along with new values.
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.81"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0991"
$ws.Range("E51").Value = "  -3.45%  "
